# Transfer BI FAST.xlsx - "fixing script automation BI fast transfer"
#
# Summary of the edit:
#  - Remove the "Rekening No BIFAST BCA" test row (old row 3); the rows
#    below it shift up.
#  - Change the account-number column (C) for the remaining "Negative"/
#    "Positive" rows from 1234567890 to 510654301.
#  - Change the nominal-transfer value in row 3 (old "BIFAST < MIN
#    Transfer" row) from 99999 to 9999.
#  - Add a new "BankTujuan" column (F) with "Digital" for every data row.
#  - Update the remembered cell selection to E13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the old row 3 ("Rekening No BIFAST BCA") -----------------
# Everything below (old rows 4-7) shifts up to become the new rows 3-6,
# carrying their existing formatting (styles, row heights) with them.
$ws.Rows("3").Delete()

# --- 2. Update the account-number (C) / nominal (D) values ------------
# Writing a brand-new literal through .Value resets the cell's style
# (Excel drops the "quotePrefix" number formatting these cells use), so
# stash the current formatting first and paste it back after the value
# change.
$ws.Range("C3").Copy()
$ws.Range("H1").PasteSpecial(-4122)

$ws.Range("D3").Copy()
$ws.Range("H2").PasteSpecial(-4122)

$ws.Range("C3").Value = 510654301
$ws.Range("C4").Value = 510654301
$ws.Range("C5").Value = 510654301
$ws.Range("C6").Value = 510654301
$ws.Range("D3").Value = 9999

$ws.Range("H1").Copy()
$ws.Range("C3:C6").PasteSpecial(-4122)

$ws.Range("H2").Copy()
$ws.Range("D3").PasteSpecial(-4122)

$ws.Range("H1:H2").Clear()

# --- 3. Add the new "BankTujuan" / "Digital" column (F) ----------------
$ws.Range("F1").Value = "BankTujuan"
$ws.Range("F2").Value = "Digital"
$ws.Range("F3").Value = "Digital"
$ws.Range("F4").Value = "Digital"
$ws.Range("F5").Value = "Digital"
$ws.Range("F6").Value = "Digital"

$ws.Columns("F").ColumnWidth = 28.8

# --- 4. Restore the saved cell selection --------------------------------
$ws.Range("E13").Select()
